$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 10.088846
$ws.Cells.Item(2, 8).Value = 30.266538
$ws.Cells.Item(2, 9).Value = 0.1151445838515654
$ws.Cells.Item(2, 10).Value = 0.1151445838515654
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 13.604331
$ws.Cells.Item(2, 14).Value = 40.812993
$ws.Cells.Item(2, 15).Value = 0.8107276168878804
$ws.Cells.Item(2, 16).Value = 0.8107276168878805
$ws.Cells.Item(2, 17).Value = 137.252000392026
$ws.Cells.Item(2, 18).Value = 1235.268003528234
$ws.Cells.Item(2, 19).Value = 0.09335089406352635
$ws.Cells.Item(2, 20).Value = 0.09335089406352637

# Row 3
$ws.Cells.Item(3, 7).Value = 10.088846
$ws.Cells.Item(3, 8).Value = 30.266538
$ws.Cells.Item(3, 9).Value = 0.1151445838515654
$ws.Cells.Item(3, 10).Value = 0.1151445838515654
$ws.Cells.Item(3, 15).Value = 0.06327311690486458
$ws.Cells.Item(3, 16).Value = 0.06327311690486459
$ws.Cells.Item(3, 17).Value = 10.711812062808
$ws.Cells.Item(3, 18).Value = 96.406308565272
$ws.Cells.Item(3, 19).Value = 0.007285556715002081
$ws.Cells.Item(3, 20).Value = 0.007285556715002082

# Row 4
$ws.Cells.Item(4, 7).Value = 10.088846
$ws.Cells.Item(4, 8).Value = 30.266538
$ws.Cells.Item(4, 9).Value = 0.1151445838515654
$ws.Cells.Item(4, 10).Value = 0.1151445838515654
$ws.Cells.Item(4, 13).Value = 1.995771333333333
$ws.Cells.Item(4, 14).Value = 5.987314
$ws.Cells.Item(4, 15).Value = 0.1189346934389115
$ws.Cells.Item(4, 16).Value = 0.1189346934389116
$ws.Cells.Item(4, 17).Value = 20.13502963321466
$ws.Cells.Item(4, 18).Value = 181.215266698932
$ws.Cells.Item(4, 19).Value = 0.01369468578153698
$ws.Cells.Item(4, 20).Value = 0.01369468578153698

# Row 5
$ws.Cells.Item(5, 7).Value = 10.088846
$ws.Cells.Item(5, 8).Value = 30.266538
$ws.Cells.Item(5, 9).Value = 0.1151445838515654
$ws.Cells.Item(5, 10).Value = 0.1151445838515654
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.1185463333333334
$ws.Cells.Item(5, 14).Value = 0.355639
$ws.Cells.Item(5, 15).Value = 0.007064572768343379
$ws.Cells.Item(5, 16).Value = 0.007064572768343379
$ws.Cells.Item(5, 17).Value = 1.195995700864667
$ws.Cells.Item(5, 18).Value = 10.763961307782
$ws.Cells.Item(5, 19).Value = 0.0008134472914999999
$ws.Cells.Item(5, 20).Value = 0.0008134472914999999

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4327250566572728
$ws.Cells.Item(6, 10).Value = 0.4327250566572729
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 13.604331
$ws.Cells.Item(6, 14).Value = 40.812993
$ws.Cells.Item(6, 15).Value = 0.8107276168878804
$ws.Cells.Item(6, 16).Value = 0.8107276168878805
$ws.Cells.Item(6, 17).Value = 515.806976405656
$ws.Cells.Item(6, 18).Value = 4642.262787650904
$ws.Cells.Item(6, 19).Value = 0.3508221539514239
$ws.Cells.Item(6, 20).Value = 0.350822153951424

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4327250566572728
$ws.Cells.Item(7, 10).Value = 0.4327250566572729
$ws.Cells.Item(7, 15).Value = 0.06327311690486458
$ws.Cells.Item(7, 16).Value = 0.06327311690486459
$ws.Cells.Item(7, 19).Value = 0.02737986309753977
$ws.Cells.Item(7, 20).Value = 0.02737986309753978

# Row 8
$ws.Cells.Item(8, 9).Value = 0.4327250566572728
$ws.Cells.Item(8, 10).Value = 0.4327250566572729
$ws.Cells.Item(8, 13).Value = 1.995771333333333
$ws.Cells.Item(8, 14).Value = 5.987314
$ws.Cells.Item(8, 15).Value = 0.1189346934389115
$ws.Cells.Item(8, 16).Value = 0.1189346934389116
$ws.Cells.Item(8, 17).Value = 75.6694891533991
$ws.Cells.Item(8, 18).Value = 681.025402380592
$ws.Cells.Item(8, 19).Value = 0.05146602195686837
$ws.Cells.Item(8, 20).Value = 0.05146602195686838

# Row 9
$ws.Cells.Item(9, 9).Value = 0.4327250566572728
$ws.Cells.Item(9, 10).Value = 0.4327250566572729
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1185463333333334
$ws.Cells.Item(9, 14).Value = 0.355639
$ws.Cells.Item(9, 15).Value = 0.007064572768343379
$ws.Cells.Item(9, 16).Value = 0.007064572768343379
$ws.Cells.Item(9, 17).Value = 4.494673480132445
$ws.Cells.Item(9, 18).Value = 40.45206132119201
$ws.Cells.Item(9, 19).Value = 0.003057017651440816
$ws.Cells.Item(9, 20).Value = 0.003057017651440816

# Row 10
$ws.Cells.Item(10, 7).Value = 15.69885766666667
$ws.Cells.Item(10, 8).Value = 47.096573
$ws.Cells.Item(10, 9).Value = 0.1791719719949428
$ws.Cells.Item(10, 10).Value = 0.1791719719949428
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 13.604331
$ws.Cells.Item(10, 14).Value = 40.812993
$ws.Cells.Item(10, 15).Value = 0.8107276168878804
$ws.Cells.Item(10, 16).Value = 0.8107276168878805
$ws.Cells.Item(10, 17).Value = 213.572456019221
$ws.Cells.Item(10, 18).Value = 1922.152104172989
$ws.Cells.Item(10, 19).Value = 0.145259665868562
$ws.Cells.Item(10, 20).Value = 0.1452596658685621

# Row 11
$ws.Cells.Item(11, 7).Value = 15.69885766666667
$ws.Cells.Item(11, 8).Value = 47.096573
$ws.Cells.Item(11, 9).Value = 0.1791719719949428
$ws.Cells.Item(11, 10).Value = 0.1791719719949428
$ws.Cells.Item(11, 15).Value = 0.06327311690486458
$ws.Cells.Item(11, 16).Value = 0.06327311690486459
$ws.Cells.Item(11, 17).Value = 16.668230729868
$ws.Cells.Item(11, 18).Value = 150.014076568812
$ws.Cells.Item(11, 19).Value = 0.01133676913011114
$ws.Cells.Item(11, 20).Value = 0.01133676913011114

# Row 12
$ws.Cells.Item(12, 7).Value = 15.69885766666667
$ws.Cells.Item(12, 8).Value = 47.096573
$ws.Cells.Item(12, 9).Value = 0.1791719719949428
$ws.Cells.Item(12, 10).Value = 0.1791719719949428
$ws.Cells.Item(12, 13).Value = 1.995771333333333
$ws.Cells.Item(12, 14).Value = 5.987314
$ws.Cells.Item(12, 15).Value = 0.1189346934389115
$ws.Cells.Item(12, 16).Value = 0.1189346934389116
$ws.Cells.Item(12, 17).Value = 31.33133009721355
$ws.Cells.Item(12, 18).Value = 281.981970874922
$ws.Cells.Item(12, 19).Value = 0.02130976356206376
$ws.Cells.Item(12, 20).Value = 0.02130976356206377

# Row 13
$ws.Cells.Item(13, 7).Value = 15.69885766666667
$ws.Cells.Item(13, 8).Value = 47.096573
$ws.Cells.Item(13, 9).Value = 0.1791719719949428
$ws.Cells.Item(13, 10).Value = 0.1791719719949428
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.1185463333333334
$ws.Cells.Item(13, 14).Value = 0.355639
$ws.Cells.Item(13, 15).Value = 0.007064572768343379
$ws.Cells.Item(13, 16).Value = 0.007064572768343379
$ws.Cells.Item(13, 17).Value = 1.861042013905223
$ws.Cells.Item(13, 18).Value = 16.749378125147
$ws.Cells.Item(13, 19).Value = 0.001265773434205856
$ws.Cells.Item(13, 20).Value = 0.001265773434205856

# Row 14
$ws.Cells.Item(14, 7).Value = 23.91632366666667
$ws.Cells.Item(14, 8).Value = 71.748971
$ws.Cells.Item(14, 9).Value = 0.2729583874962189
$ws.Cells.Item(14, 10).Value = 0.2729583874962189
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 13.604331
$ws.Cells.Item(14, 14).Value = 40.812993
$ws.Cells.Item(14, 15).Value = 0.8107276168878804
$ws.Cells.Item(14, 16).Value = 0.8107276168878805
$ws.Cells.Item(14, 17).Value = 325.365583464467
$ws.Cells.Item(14, 18).Value = 2928.290251180203
$ws.Cells.Item(14, 19).Value = 0.2212949030043682
$ws.Cells.Item(14, 20).Value = 0.2212949030043682

# Row 15
$ws.Cells.Item(15, 7).Value = 23.91632366666667
$ws.Cells.Item(15, 8).Value = 71.748971
$ws.Cells.Item(15, 9).Value = 0.2729583874962189
$ws.Cells.Item(15, 10).Value = 0.2729583874962189
$ws.Cells.Item(15, 15).Value = 0.06327311690486458
$ws.Cells.Item(15, 16).Value = 0.06327311690486459
$ws.Cells.Item(15, 17).Value = 25.393108820436
$ws.Cells.Item(15, 18).Value = 228.537979383924
$ws.Cells.Item(15, 19).Value = 0.01727092796221159
$ws.Cells.Item(15, 20).Value = 0.01727092796221159

# Row 16
$ws.Cells.Item(16, 7).Value = 23.91632366666667
$ws.Cells.Item(16, 8).Value = 71.748971
$ws.Cells.Item(16, 9).Value = 0.2729583874962189
$ws.Cells.Item(16, 10).Value = 0.2729583874962189
$ws.Cells.Item(16, 13).Value = 1.995771333333333
$ws.Cells.Item(16, 14).Value = 5.987314
$ws.Cells.Item(16, 15).Value = 0.1189346934389115
$ws.Cells.Item(16, 16).Value = 0.1189346934389116
$ws.Cells.Item(16, 17).Value = 47.73151317265489
$ws.Cells.Item(16, 18).Value = 429.5836185538939
$ws.Cells.Item(16, 19).Value = 0.03246422213844242
$ws.Cells.Item(16, 20).Value = 0.03246422213844243

# Row 17
$ws.Cells.Item(17, 7).Value = 23.91632366666667
$ws.Cells.Item(17, 8).Value = 71.748971
$ws.Cells.Item(17, 9).Value = 0.2729583874962189
$ws.Cells.Item(17, 10).Value = 0.2729583874962189
$ws.Cells.Item(17, 11).Value = 2
$ws.Cells.Item(17, 12).Value = 0.6666666666666666
$ws.Cells.Item(17, 13).Value = 0.1185463333333334
$ws.Cells.Item(17, 14).Value = 0.355639
$ws.Cells.Item(17, 15).Value = 0.007064572768343379
$ws.Cells.Item(17, 16).Value = 0.007064572768343379
$ws.Cells.Item(17, 17).Value = 2.835192477496556
$ws.Cells.Item(17, 18).Value = 25.516732297469
$ws.Cells.Item(17, 19).Value = 0.001928334391196708
$ws.Cells.Item(17, 20).Value = 0.001928334391196708
